$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 4 had their D (Fecha), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) values swapped.

$cols = @("D", "N", "O", "P", "S")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow4 = $ws.Range($col + "4")

    $valRow2 = $cellRow2.Value2
    $valRow4 = $cellRow4.Value2

    $cellRow2.Value2 = $valRow4
    $cellRow4.Value2 = $valRow2
}
